$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 34.446933
$ws.Range("H2").Value = 103.340799
$ws.Range("I2").Value = 0.3406054910909001
$ws.Range("J2").Value = 0.3406054910909001
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 112.708133
$ws.Range("N2").Value = 338.124399
$ws.Range("O2").Value = 0.2121524692929861
$ws.Range("P2").Value = 0.2121524692929861
$ws.Range("Q2").Value = 3882.44950600609
$ws.Range("R2").Value = 34942.04555405481
$ws.Range("S2").Value = 0.07226029598968463
$ws.Range("T2").Value = 0.07226029598968463

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 34.446933
$ws.Range("H3").Value = 103.340799
$ws.Range("I3").Value = 0.3406054910909001
$ws.Range("J3").Value = 0.3406054910909001
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 189.57842
$ws.Range("N3").Value = 568.7352599999999
$ws.Range("O3").Value = 0.3568467408440064
$ws.Range("P3").Value = 0.3568467408440064
$ws.Range("Q3").Value = 6530.395131985859
$ws.Range("R3").Value = 58773.55618787273
$ws.Range("S3").Value = 0.12154395940936
$ws.Range("T3").Value = 0.12154395940936

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 34.446933
$ws.Range("H4").Value = 103.340799
$ws.Range("I4").Value = 0.3406054910909001
$ws.Range("J4").Value = 0.3406054910909001
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 110.6512236666667
$ws.Range("N4").Value = 331.953671
$ws.Range("O4").Value = 0.2082807132576123
$ws.Range("P4").Value = 0.2082807132576123
$ws.Range("Q4").Value = 3811.595288013681
$ws.Range("R4").Value = 34304.35759212313
$ws.Range("S4").Value = 0.07094155462387199
$ws.Range("T4").Value = 0.07094155462387197

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 34.446933
$ws.Range("H5").Value = 103.340799
$ws.Range("I5").Value = 0.3406054910909001
$ws.Range("J5").Value = 0.3406054910909001
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 118.3222806666667
$ws.Range("N5").Value = 354.966842
$ws.Range("O5").Value = 0.2227200766053953
$ws.Range("P5").Value = 0.2227200766053952
$ws.Range("Q5").Value = 4075.839674531862
$ws.Range("R5").Value = 36682.55707078676
$ws.Range("S5").Value = 0.07585968106798353
$ws.Range("T5").Value = 0.07585968106798353

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 40.11196
$ws.Range("H6").Value = 120.33588
$ws.Range("I6").Value = 0.396620327110647
$ws.Range("J6").Value = 0.396620327110647
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 112.708133
$ws.Range("N6").Value = 338.124399
$ws.Range("O6").Value = 0.2121524692929861
$ws.Range("P6").Value = 0.2121524692929861
$ws.Range("Q6").Value = 4520.944122570681
$ws.Range("R6").Value = 40688.49710313613
$ws.Range("S6").Value = 0.08414398176831564
$ws.Range("T6").Value = 0.08414398176831563

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 40.11196
$ws.Range("H7").Value = 120.33588
$ws.Range("I7").Value = 0.396620327110647
$ws.Range("J7").Value = 0.396620327110647
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 189.57842
$ws.Range("N7").Value = 568.7352599999999
$ws.Range("O7").Value = 0.3568467408440064
$ws.Range("P7").Value = 0.3568467408440064
$ws.Range("Q7").Value = 7604.361999903199
$ws.Range("R7").Value = 68439.2579991288
$ws.Range("S7").Value = 0.1415326710819181
$ws.Range("T7").Value = 0.1415326710819181

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 40.11196
$ws.Range("H8").Value = 120.33588
$ws.Range("I8").Value = 0.396620327110647
$ws.Range("J8").Value = 0.396620327110647
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 110.6512236666667
$ws.Range("N8").Value = 331.953671
$ws.Range("O8").Value = 0.2082807132576123
$ws.Range("P8").Value = 0.2082807132576123
$ws.Range("Q8").Value = 4438.437457668387
$ws.Range("R8").Value = 39945.93711901548
$ws.Range("S8").Value = 0.08260836462307307
$ws.Range("T8").Value = 0.08260836462307304

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 40.11196
$ws.Range("H9").Value = 120.33588
$ws.Range("I9").Value = 0.396620327110647
$ws.Range("J9").Value = 0.396620327110647
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 118.3222806666667
$ws.Range("N9").Value = 354.966842
$ws.Range("O9").Value = 0.2227200766053953
$ws.Range("P9").Value = 0.2227200766053952
$ws.Range("Q9").Value = 4746.138589210107
$ws.Range("R9").Value = 42715.24730289096
$ws.Range("S9").Value = 0.08833530963734024
$ws.Range("T9").Value = 0.08833530963734021

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.3464193333333334
$ws.Range("H10").Value = 1.039258
$ws.Range("I10").Value = 0.003425336216532898
$ws.Range("J10").Value = 0.003425336216532898
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 112.708133
$ws.Range("N10").Value = 338.124399
$ws.Range("O10").Value = 0.2121524692929861
$ws.Range("P10").Value = 0.2121524692929861
$ws.Range("Q10").Value = 39.04427629510467
$ws.Range("R10").Value = 351.3984866559421
$ws.Range("S10").Value = 0.0007266935364961488
$ws.Range("T10").Value = 0.0007266935364961487

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.3464193333333334
$ws.Range("H11").Value = 1.039258
$ws.Range("I11").Value = 0.003425336216532898
$ws.Range("J11").Value = 0.003425336216532898
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 189.57842
$ws.Range("N11").Value = 568.7352599999999
$ws.Range("O11").Value = 0.3568467408440064
$ws.Range("P11").Value = 0.3568467408440064
$ws.Range("Q11").Value = 65.67362987078666
$ws.Range("R11").Value = 591.06266883708
$ws.Range("S11").Value = 0.001222320065164704
$ws.Range("T11").Value = 0.001222320065164704

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.3464193333333334
$ws.Range("H12").Value = 1.039258
$ws.Range("I12").Value = 0.003425336216532898
$ws.Range("J12").Value = 0.003425336216532898
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 110.6512236666667
$ws.Range("N12").Value = 331.953671
$ws.Range("O12").Value = 0.2082807132576123
$ws.Range("P12").Value = 0.2082807132576123
$ws.Range("Q12").Value = 38.33172313512422
$ws.Range("R12").Value = 344.985508216118
$ws.Range("S12").Value = 0.000713431470326603
$ws.Range("T12").Value = 0.000713431470326603

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.3464193333333334
$ws.Range("H13").Value = 1.039258
$ws.Range("I13").Value = 0.003425336216532898
$ws.Range("J13").Value = 0.003425336216532898
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 118.3222806666667
$ws.Range("N13").Value = 354.966842
$ws.Range("O13").Value = 0.2227200766053953
$ws.Range("P13").Value = 0.2227200766053952
$ws.Range("Q13").Value = 40.98912558702622
$ws.Range("R13").Value = 368.902130283236
$ws.Range("S13").Value = 0.0007628911445454417
$ws.Range("T13").Value = 0.0007628911445454416

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 26.22909066666666
$ws.Range("H14").Value = 78.687272
$ws.Range("I14").Value = 0.25934884558192
$ws.Range("J14").Value = 0.25934884558192
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 112.708133
$ws.Range("N14").Value = 338.124399
$ws.Range("O14").Value = 0.2121524692929861
$ws.Range("P14").Value = 0.2121524692929861
$ws.Range("Q14").Value = 2956.231839327726
$ws.Range("R14").Value = 26606.08655394953
$ws.Range("S14").Value = 0.05502149799848966
$ws.Range("T14").Value = 0.05502149799848965

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 26.22909066666666
$ws.Range("H15").Value = 78.687272
$ws.Range("I15").Value = 0.25934884558192
$ws.Range("J15").Value = 0.25934884558192
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 189.57842
$ws.Range("N15").Value = 568.7352599999999
$ws.Range("O15").Value = 0.3568467408440064
$ws.Range("P15").Value = 0.3568467408440064
$ws.Range("Q15").Value = 4972.469566623412
$ws.Range("R15").Value = 44752.22609961071
$ws.Range("S15").Value = 0.09254779028756363
$ws.Range("T15").Value = 0.09254779028756363

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 26.22909066666666
$ws.Range("H16").Value = 78.687272
$ws.Range("I16").Value = 0.25934884558192
$ws.Range("J16").Value = 0.25934884558192
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 110.6512236666667
$ws.Range("N16").Value = 331.953671
$ws.Range("O16").Value = 0.2082807132576123
$ws.Range("P16").Value = 0.2082807132576123
$ws.Range("Q16").Value = 2902.280977930612
$ws.Range("R16").Value = 26120.52880137551
$ws.Range("S16").Value = 0.05401736254034064
$ws.Range("T16").Value = 0.05401736254034063

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 26.22909066666666
$ws.Range("H17").Value = 78.687272
$ws.Range("I17").Value = 0.25934884558192
$ws.Range("J17").Value = 0.25934884558192
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 118.3222806666667
$ws.Range("N17").Value = 354.966842
$ws.Range("O17").Value = 0.2227200766053953
$ws.Range("P17").Value = 0.2227200766053952
$ws.Range("Q17").Value = 3103.48582749278
$ws.Range("R17").Value = 27931.37244743502
$ws.Range("S17").Value = 0.05776219475552603
$ws.Range("T17").Value = 0.05776219475552603
